# Auto-generated edit script applying the Diabolos_Profits.xlsx diff
$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# ---- Sheet index 1 ----
$ws = $sheets.Item(1)
$ws.Range("H5").Value = 1250.5
$ws.Range("I5").Value = 1001
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 1001
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = -886
$ws.Range("N5").Value = -1730
$ws.Range("H32").Value = 6039.5713
$ws.Range("I32").Value = 2990
$ws.Range("J32").Value = 6547.8335
$ws.Range("K32").Value = 2990
$ws.Range("L32").Value = 6547.8335
$ws.Range("M32").Value = -2664
$ws.Range("N32").Value = -7199.8335
$ws.Range("H100").Value = 9000
$ws.Range("I100").Value = 9000
$ws.Range("K100").Value = 9000
$ws.Range("M100").Value = -8459
$ws.Range("H106").Value = 7895.5
$ws.Range("I106").Value = 3685.2727
$ws.Range("K106").Value = 3685.2727
$ws.Range("M106").Value = -3054.2727
$ws.Range("H111").Value = 36853.555
$ws.Range("J111").Value = 59203.25
$ws.Range("L111").Value = 177609.75
$ws.Range("N111").Value = -183743.75
$ws.Range("H127").Value = 1335.8235
$ws.Range("I127").Value = 1108.1428
$ws.Range("J127").Value = 2398.3333
$ws.Range("K127").Value = 3324.4284
$ws.Range("L127").Value = 7194.999899999999
$ws.Range("M127").Value = 1635.5716
$ws.Range("N127").Value = -17114.9999
$ws.Range("H131").Value = 13853.571
$ws.Range("I131").Value = 5750
$ws.Range("J131").Value = 14476.923
$ws.Range("K131").Value = 17250
$ws.Range("L131").Value = 43430.769
$ws.Range("M131").Value = -12210
$ws.Range("N131").Value = -53510.769
$ws.Range("H132").Value = 4553
$ws.Range("I132").Value = 5522.4
$ws.Range("K132").Value = 16567.2
$ws.Range("M132").Value = -14037.2
$ws.Range("H141").Value = 1589.6842
$ws.Range("I141").Value = 680.26666
$ws.Range("K141").Value = 2040.79998
$ws.Range("M141").Value = 3139.20002

# ---- Sheet index 2 ----
$ws = $sheets.Item(2)
$ws.Range("H32").Value = 187387.19
$ws.Range("I32").Value = 210691.44
$ws.Range("J32").Value = 953.1667
$ws.Range("K32").Value = 210691.44
$ws.Range("L32").Value = 953.1667
$ws.Range("M32").Value = -210404.44
$ws.Range("N32").Value = -1527.1667
$ws.Range("H61").Value = 2699.818
$ws.Range("J61").Value = 3751.3333
$ws.Range("L61").Value = 3751.3333
$ws.Range("N61").Value = -4175.3333
$ws.Range("H110").Value = 71440750
$ws.Range("I110").Value = 111112510
$ws.Range("J110").Value = 31591.8
$ws.Range("K110").Value = 111112510
$ws.Range("L110").Value = 31591.8
$ws.Range("M110").Value = -111110465
$ws.Range("N110").Value = -35681.8
$ws.Range("H136").Value = 2699.818
$ws.Range("J136").Value = 3751.3333
$ws.Range("L136").Value = 11253.9999
$ws.Range("N136").Value = -16353.9999
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet index 3 ----
$ws = $sheets.Item(3)
$ws.Range("H86").Value = 2271
$ws.Range("I86").Value = 2299.5715
$ws.Range("K86").Value = 2299.5715
$ws.Range("M86").Value = -1176.5715
$ws.Range("H89").Value = 2271
$ws.Range("I89").Value = 2299.5715
$ws.Range("K89").Value = 11497.8575
$ws.Range("M89").Value = -5881.8575

# ---- Sheet index 4 ----
$ws = $sheets.Item(4)
$ws.Range("H31").Value = 2568.6333
$ws.Range("I31").Value = 1715.25
$ws.Range("J31").Value = 3315.3438
$ws.Range("K31").Value = 1715.25
$ws.Range("L31").Value = 3315.3438
$ws.Range("M31").Value = -1420.25
$ws.Range("N31").Value = -3905.3438
$ws.Range("H34").Value = 2568.6333
$ws.Range("I34").Value = 1715.25
$ws.Range("J34").Value = 3315.3438
$ws.Range("K34").Value = 1715.25
$ws.Range("L34").Value = 3315.3438
$ws.Range("M34").Value = -1513.25
$ws.Range("N34").Value = -3719.3438
$ws.Range("H99").Value = 2249.1904
$ws.Range("I99").Value = 1958.75
$ws.Range("J99").Value = 2636.4443
$ws.Range("K99").Value = 1958.75
$ws.Range("L99").Value = 2636.4443
$ws.Range("M99").Value = -460.75
$ws.Range("N99").Value = -5632.4443
$ws.Range("H107").Value = 639.0769
$ws.Range("J107").Value = 475.33334
$ws.Range("L107").Value = 475.33334
$ws.Range("N107").Value = -4315.33334
$ws.Range("H126").Value = 2249.1904
$ws.Range("I126").Value = 1958.75
$ws.Range("J126").Value = 2636.4443
$ws.Range("K126").Value = 5876.25
$ws.Range("L126").Value = 7909.3329
$ws.Range("M126").Value = -3406.25
$ws.Range("N126").Value = -12849.3329
$ws.Range("H132").Value = 3749.92
$ws.Range("I132").Value = 2330.1333
$ws.Range("J132").Value = 5879.6
$ws.Range("K132").Value = 6990.3999
$ws.Range("L132").Value = 17638.8
$ws.Range("M132").Value = -4460.3999
$ws.Range("N132").Value = -22698.8

# ---- Sheet index 5 ----
$ws = $sheets.Item(5)
$ws.Range("H92").Value = 1428.375
$ws.Range("J92").Value = 1662.8334
$ws.Range("L92").Value = 4988.5002
$ws.Range("N92").Value = -7484.5002
$ws.Range("H105").Value = 15718.956
$ws.Range("I105").Value = 15063
$ws.Range("K105").Value = 45189
$ws.Range("M105").Value = -42568
$ws.Range("H107").Value = 1376.6923
$ws.Range("I107").Value = 1376.6923
$ws.Range("K107").Value = 4130.0769
$ws.Range("M107").Value = -2210.0769
$ws.Range("H113").Value = 1464.3846
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 1553.0834
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 4659.2502
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -8999.2502

# ---- Sheet index 6 ----
$ws = $sheets.Item(6)
$ws.Range("H70").Value = 7907.9614
$ws.Range("I70").Value = 8182.0625
$ws.Range("K70").Value = 8182.0625
$ws.Range("M70").Value = -7912.0625
$ws.Range("H73").Value = 7907.9614
$ws.Range("I73").Value = 8182.0625
$ws.Range("K73").Value = 8182.0625
$ws.Range("M73").Value = -7246.0625
$ws.Range("H102").Value = 2501.7727
$ws.Range("I102").Value = 1386.6666
$ws.Range("J102").Value = 3839.9
$ws.Range("K102").Value = 1386.6666
$ws.Range("L102").Value = 3839.9
$ws.Range("M102").Value = 235.3334
$ws.Range("N102").Value = -7083.9
$ws.Range("H122").Value = 26318392
$ws.Range("I122").Value = 2348.6428
$ws.Range("K122").Value = 7045.928400000001
$ws.Range("M122").Value = -4595.928400000001

# ---- Sheet index 7 ----
$ws = $sheets.Item(7)
$ws.Range("H7").Value = 29415200
$ws.Range("I7").Value = 71431100
$ws.Range("J7").Value = 4065.2
$ws.Range("K7").Value = 71431100
$ws.Range("L7").Value = 4065.2
$ws.Range("M7").Value = -71430988
$ws.Range("N7").Value = -4289.2
$ws.Range("H16").Value = 3210.2104
$ws.Range("I16").Value = 2882.9333
$ws.Range("J16").Value = 4437.5
$ws.Range("K16").Value = 2882.9333
$ws.Range("L16").Value = 4437.5
$ws.Range("M16").Value = -2712.9333
$ws.Range("N16").Value = -4777.5
$ws.Range("H22").Value = 1614183.5
$ws.Range("I22").Value = 1306
$ws.Range("J22").Value = 4609527.5
$ws.Range("K22").Value = 1306
$ws.Range("L22").Value = 4609527.5
$ws.Range("M22").Value = -1011
$ws.Range("N22").Value = -4610117.5
$ws.Range("H27").Value = 1614183.5
$ws.Range("I27").Value = 1306
$ws.Range("J27").Value = 4609527.5
$ws.Range("K27").Value = 1306
$ws.Range("L27").Value = 4609527.5
$ws.Range("M27").Value = -1199
$ws.Range("N27").Value = -4609741.5
$ws.Range("H68").Value = 8162.5
$ws.Range("J68").Value = 10328.917
$ws.Range("L68").Value = 10328.917
$ws.Range("N68").Value = -11826.917
$ws.Range("H71").Value = 8162.5
$ws.Range("J71").Value = 10328.917
$ws.Range("L71").Value = 51644.585
$ws.Range("N71").Value = -59132.585
$ws.Range("H93").Value = 1759.9286
$ws.Range("I93").Value = 1852.6666
$ws.Range("K93").Value = 1852.6666
$ws.Range("M93").Value = -604.6666
$ws.Range("H122").Value = 3506.093
$ws.Range("I122").Value = 2731.5925
$ws.Range("J122").Value = 4813.0625
$ws.Range("K122").Value = 8194.7775
$ws.Range("L122").Value = 14439.1875
$ws.Range("M122").Value = -5744.7775
$ws.Range("N122").Value = -19339.1875
$ws.Range("H126").Value = 29415200
$ws.Range("I126").Value = 71431100
$ws.Range("J126").Value = 4065.2
$ws.Range("K126").Value = 214293300
$ws.Range("L126").Value = 12195.6
$ws.Range("M126").Value = -214290830
$ws.Range("N126").Value = -17135.6

# ---- Sheet index 8 ----
$ws = $sheets.Item(8)
$ws.Range("H96").Value = 5208.5713
$ws.Range("J96").Value = 5434.222
$ws.Range("L96").Value = 5434.222
$ws.Range("N96").Value = -8180.222
$ws.Range("H100").Value = 406.81818
$ws.Range("I100").Value = 392.6
$ws.Range("J100").Value = 549
$ws.Range("K100").Value = 785.2
$ws.Range("L100").Value = 1098
$ws.Range("M100").Value = -244.2
$ws.Range("N100").Value = -2180
$ws.Range("H126").Value = 2097.5454
$ws.Range("I126").Value = 1997.7142
$ws.Range("K126").Value = 5993.142599999999
$ws.Range("M126").Value = -3523.142599999999
$ws.Range("H136").Value = 2801.4814
$ws.Range("I136").Value = 1404.1818
$ws.Range("J136").Value = 8949.6
$ws.Range("K136").Value = 4212.5454
$ws.Range("L136").Value = 26848.8
$ws.Range("M136").Value = -1662.5454
$ws.Range("N136").Value = -31948.8
